$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column retains its text formatting so values such as
# "1.00" or "9.40" are not silently converted to numbers, matching the
# original inline-string cell content.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '71.043.88'
$ws.Range("E2").Value = '  +1.93%  '
$ws.Range("D3").Value = '3.584.28'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '585.09'
$ws.Range("E5").Value = '  +2.28%  '
$ws.Range("D6").Value = '186.40'
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("D7").Value = '3.573.08'
$ws.Range("E7").Value = '  +1.27%  '
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  +19.51%  '
$ws.Range("D11").Value = '0.652'
$ws.Range("E11").Value = '  +1.20%  '
$ws.Range("D12").Value = '54.36'
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("E13").Value = '  +7.75%  '
$ws.Range("D14").Value = '9.54'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").Value = '4.147.88'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").Value = '70.933.41'
$ws.Range("E16").Value = '  +1.90%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").Value = '19.30'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.597.43'
$ws.Range("E18").Value = '  +1.81%  '
$ws.Range("D19").Value = '573.08'
$ws.Range("E19").Value = '  +12.64%  '
$ws.Range("D20").Value = '12.37'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").Value = '0.121'
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("E22").Value = '  -2.29%  '
$ws.Range("D23").Value = '17.46'
$ws.Range("E23").Value = '  -12.41%  '
$ws.Range("D24").Value = '5.09'
$ws.Range("E24").Value = '  +3.27%  '
$ws.Range("E25").Value = '  +5.60%  '
$ws.Range("D26").Value = '95.01'
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("D27").Value = '11.33'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("D28").Value = '2.93'
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").Value = '9.11'
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("D30").Value = '32.33'
$ws.Range("E30").Value = '  +2.80%  '
$ws.Range("D31").Value = '7.22'
$ws.Range("E31").Value = '  -4.06%  '
$ws.Range("D32").Value = '12.28'
$ws.Range("E32").Value = '  -1.37%  '
$ws.Range("D33").Value = '0.115'
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("D34").Value = '64.09'
$ws.Range("E34").Value = '  -2.11%  '
$ws.Range("D35").Value = '3.33'
$ws.Range("E35").Value = '  +5.45%  '
$ws.Range("D36").Value = '554.98'
$ws.Range("E36").Value = '  -2.29%  '
$ws.Range("D37").Value = '0.414'
$ws.Range("E37").Value = '  +2.99%  '
$ws.Range("E38").Value = '  +4.91%  '
$ws.Range("D39").Value = '37.52'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("D41").Value = '3.501.51'
$ws.Range("E41").Value = '  +10.55%  '
$ws.Range("D42").Value = '3.20'
$ws.Range("E42").Value = '  +1.92%  '
$ws.Range("D43").Value = '3.44'
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("D46").Value = '0.0446'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("D48").Value = '9.40'
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("E49").Value = '  +2.27%  '
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("E51").Value = '  -1.18%  '

$wb.Save()
